# Fill row 1 (A1:K1) with the new data values, then leave the selection
# on the next empty cell (L1) — matching the state Excel leaves behind
# after typing values across a row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1, 1, 1, 11, 1, 1, 1, 1, 1, 1, 1)
for ($col = 1; $col -le $values.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $values[$col - 1]
}

$ws.Range("L1").Select()
